$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E5").Value = 47752663.80848
$ws.Range("E6").Value = 142160477.54952
$ws.Range("E7").Value = 47752638.67548
$ws.Range("E8").Value = 81822994.8
